$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 893
$ws1.Range("F3").Value = 4527

# Update the same values on the "全部类型" sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 893
$ws4.Range("F3").Value = 4527
